$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings are stored as text (matches the
# original inline-string cells), not coerced into floating point numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.768.24"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "3.090.03"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "540.70"
$ws.Range("E5").Value = "  -2.90%  "
$ws.Range("D6").Value = "136.91"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.083.75"
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("D11").Value = "6.26"
$ws.Range("E11").Value = "  -5.90%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "0.0000228"
$ws.Range("E13").Value = "  +4.19%  "
$ws.Range("D14").Value = "34.86"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("D15").Value = "3.587.51"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").Value = "63.753.69"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "3.088.35"
$ws.Range("E18").Value = "  -0.36%  "
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "489.23"
$ws.Range("E20").Value = "  -4.06%  "
$ws.Range("D21").Value = "13.50"
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("D22").Value = "0.703"
$ws.Range("E22").Value = "  -1.18%  "
$ws.Range("D23").Value = "7.22"
$ws.Range("E23").Value = "  -1.30%  "
$ws.Range("D24").Value = "79.87"
$ws.Range("E24").Value = "  +2.16%  "
$ws.Range("D25").Value = "12.29"
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("D28").Value = "8.34"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").Value = "26.30"
$ws.Range("E30").Value = "  -0.56%  "
$ws.Range("E31").Value = "  -2.93%  "
$ws.Range("D32").Value = "1.13"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("E33").Value = "  -4.87%  "
$ws.Range("D34").Value = "57.11"
$ws.Range("E34").Value = "  -4.50%  "
$ws.Range("D35").Value = "5.50"
$ws.Range("E35").Value = "  +4.67%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D36").Value = "496.46"
$ws.Range("E36").Value = "  -7.16%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").Value = "6.08"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").Value = "3.290.22"
$ws.Range("E38").Value = "  +6.97%  "
$ws.Range("E39").Value = "  -3.67%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("E41").Value = "  -2.69%  "
$ws.Range("D42").Value = "8.17"
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("E43").Value = "  -2.33%  "
$ws.Range("D44").Value = "0.259"
$ws.Range("E44").Value = "  +0.78%  "
$ws.Range("D46").Value = "2.09"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").Value = "0.0₃0540"
$ws.Range("E47").Value = "  +4.88%  "
$ws.Range("D48").Value = "24.96"
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("D49").Value = "121.79"
$ws.Range("E49").Value = "  -0.94%  "
$ws.Range("E50").Value = "  +2.01%  "
$ws.Range("D51").Value = "2.36"
$ws.Range("E51").Value = "  -2.63%  "
